$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (19) - "Sub brand" - shifting
# everything from S onward one column to the right.
$ws.Columns.Item(19).Insert()

# New column header, matching the style of the surrounding header row (S1).
$ws.Cells.Item(1, 19).Value = "Sub brand"
$ws.Cells.Item(1, 19).Font.Bold = $true
$ws.Cells.Item(1, 19).HorizontalAlignment = -4108
$ws.Cells.Item(1, 19).VerticalAlignment = -4108

# Update the view: scroll so column H is the first visible column, and
# select S1 (the new header cell).
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("S1").Select()
